# Insert a new weekly price record for "Feria Lagunitas de Puerto Montt - Zapallo italiano"
# right after the existing row for A71 (shifts all subsequent rows down by one),
# then populate the new row's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a whole new row at position 71; everything currently at/after row 71
# (through row 166) shifts down to rows 72-167, and the sheet's used range
# grows from A1:R166 to A1:R167 automatically.
$ws.Rows("71:71").Insert()

# Preserve the date-formatted style used throughout column D.
$ws.Range("D71").NumberFormat = $ws.Range("D72").NumberFormat

# Populate the newly inserted row with the new observation.
$ws.Range("A71").Value2 = 4
$ws.Range("B71").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C71").Value2 = "Los Lagos"
$ws.Range("D71").Value2 = 44546
$ws.Range("E71").Value2 = 10
$ws.Range("F71").Value2 = 100112032
$ws.Range("G71").Value2 = "Zapallo italiano"
$ws.Range("H71").Value2 = "Sin especificar"
$ws.Range("I71").Value2 = "Primera"
$ws.Range("J71").Value2 = 120
$ws.Range("K71").Value2 = 11000
$ws.Range("L71").Value2 = 12000
$ws.Range("M71").Value2 = 11500
$ws.Range("N71").Value2 = "`$/caja 50 unidades"
$ws.Range("O71").Value2 = "Región de O'Higgins"
$ws.Range("P71").Value2 = 230
$ws.Range("Q71").Value2 = 50
$ws.Range("R71").Value2 = "Hortaliza"
